$d = $word.ActiveDocument

# 1. Summary paragraph: "60%" -> "59%" of all emails that are L1
$d.Content.Find.Execute(
    "The classifier predicts L1 on 60% of all emails",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The classifier predicts L1 on 59% of all emails", 2) | Out-Null

# 2. Model description paragraph: replace chi2 feature-selection clause
#    with "using only 10,000 features"
$d.Content.Find.Execute(
    "regular expression stemmer, and chi2 feature selection for selecting values with p < 0.1.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "regular expression stemmer, and using only 10,000 features.", 2) | Out-Null

# 3. Move the "_GoBack" bookmark from its old spot (an empty paragraph
#    right before the "Background:" heading) to a collapsed point
#    between "53" and "% of emails Auto-triaged" further down the doc.
$d.Bookmarks.ShowHidden = $true

$target = $d.Content
$target.Find.Execute(
    "% of emails Auto-triaged", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null

$bookmarkRange = $d.Range($target.Start, $target.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
